$wb = $excel.ActiveWorkbook

# ALC row 92 (Leve Item ID 19901)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 1755.8
$ws.Range("I92").Value = 1395.3334
$ws.Range("J92").Value = 5000
$ws.Range("K92").Value = 1395.3334
$ws.Range("L92").Value = 5000
$ws.Range("M92").Value = -147.3334
$ws.Range("N92").Value = -7496

# ALC row 100 (Leve Item ID 19906)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 40002300
$ws.Range("I100").Value = 40002300
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 40002300
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -40001759
$ws.Range("N100").ClearContents()

# ARM row 97 (Leve Item ID 19941)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 2115.7273
$ws.Range("I97").Value = 1472
$ws.Range("J97").Value = 2652.1667
$ws.Range("K97").Value = 1472
$ws.Range("L97").Value = 2652.1667
$ws.Range("M97").Value = -976
$ws.Range("N97").Value = -3644.1667

# ARM row 102 (Leve Item ID 19945)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").ClearContents()

# ARM row 115 (Leve Item ID 27104)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H115").Value = 29908.334
$ws.Range("J115").Value = 29908.334
$ws.Range("L115").Value = 29908.334
$ws.Range("N115").Value = -33042.334

# ARM row 122 (Leve Item ID 36168)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 6376.8237
$ws.Range("I122").Value = 5337.846
$ws.Range("J122").Value = 9753.5
$ws.Range("K122").Value = 16013.538
$ws.Range("L122").Value = 29260.5
$ws.Range("M122").Value = -13563.538
$ws.Range("N122").Value = -34160.5

# ARM row 132 (Leve Item ID 43997)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 3634.3333
$ws.Range("I132").Value = 1633.3334
$ws.Range("J132").Value = 4301.3335
$ws.Range("K132").Value = 4900.0002
$ws.Range("L132").Value = 12904.0005
$ws.Range("M132").Value = -2370.0002
$ws.Range("N132").Value = -17964.0005

# BSM row 94 (Leve Item ID 19939)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 689.875
$ws.Range("I94").Value = 694.8333
$ws.Range("J94").Value = 675
$ws.Range("K94").Value = 694.8333
$ws.Range("L94").Value = 675
$ws.Range("M94").Value = -243.8333
$ws.Range("N94").Value = -1577

# BSM row 99 (Leve Item ID 19943)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 5951.6665
$ws.Range("I99").Value = 3010
$ws.Range("J99").Value = 6540
$ws.Range("K99").Value = 3010
$ws.Range("L99").Value = 6540
$ws.Range("M99").Value = -1512
$ws.Range("N99").Value = -9536

# BSM row 123 (Leve Item ID 35320)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H123").Value = 30296.666
$ws.Range("J123").Value = 30296.666
$ws.Range("L123").Value = 30296.666
$ws.Range("N123").Value = -40096.666

# CRP row 58 (Leve Item ID 44021)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2080.1785
$ws.Range("I58").Value = 1569.8
$ws.Range("J58").Value = 6333.3335
$ws.Range("K58").Value = 1569.8
$ws.Range("L58").Value = 6333.3335
$ws.Range("M58").Value = -1366.8
$ws.Range("N58").Value = -6739.3335

# CRP row 132 (Leve Item ID 44019)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2560.48
$ws.Range("I132").Value = 1609.95
$ws.Range("J132").Value = 6362.6
$ws.Range("K132").Value = 4829.85
$ws.Range("L132").Value = 19087.8
$ws.Range("M132").Value = -2299.85
$ws.Range("N132").Value = -24147.8

# CRP row 134 (Leve Item ID 44020)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 6100.88
$ws.Range("I134").Value = 7271.294
$ws.Range("K134").Value = 21813.882
$ws.Range("M134").Value = -19278.882

# CRP row 136 (Leve Item ID 44021)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 2080.1785
$ws.Range("I136").Value = 1569.8
$ws.Range("J136").Value = 6333.3335
$ws.Range("K136").Value = 4709.4
$ws.Range("L136").Value = 19000.0005
$ws.Range("M136").Value = -2159.4
$ws.Range("N136").Value = -24100.0005

# CUL row 14 (Leve Item ID 12886)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 116.90909
$ws.Range("I14").Value = 116.90909
$ws.Range("K14").Value = 350.72727
$ws.Range("M14").Value = -177.72727

# CUL row 38 (Leve Item ID 4860)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 119.64706
$ws.Range("I38").Value = 59.333332
$ws.Range("J38").Value = 187.5
$ws.Range("K38").Value = 177.999996
$ws.Range("L38").Value = 562.5
$ws.Range("M38").Value = 169.000004
$ws.Range("N38").Value = -1256.5

# CUL row 98 (Leve Item ID 19843)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 195.25
$ws.Range("I98").Value = 183.72728
$ws.Range("J98").Value = 220.6
$ws.Range("K98").Value = 551.18184
$ws.Range("L98").Value = 661.8
$ws.Range("M98").Value = 946.81816
$ws.Range("N98").Value = -3657.8

# CUL row 112 (Leve Item ID 27855)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H112").Value = 2441.125
$ws.Range("I112").Value = 975
$ws.Range("J112").Value = 2929.8333
$ws.Range("K112").Value = 2925
$ws.Range("L112").Value = 8789.499899999999
$ws.Range("M112").Value = -1817
$ws.Range("N112").Value = -11005.4999

# GSM row 74 (Leve Item ID 10972)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H74").Value = 39766.668
$ws.Range("J74").Value = 39766.668
$ws.Range("L74").Value = 39766.668
$ws.Range("N74").Value = -41638.668

# GSM row 77 (Leve Item ID 10972)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H77").Value = 39766.668
$ws.Range("J77").Value = 39766.668
$ws.Range("L77").Value = 119300.004
$ws.Range("N77").Value = -128660.004

# GSM row 113 (Leve Item ID 27710)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 2621.1
$ws.Range("I113").Value = 2630.1428
$ws.Range("J113").Value = 2600
$ws.Range("K113").Value = 2630.1428
$ws.Range("L113").Value = 2600
$ws.Range("M113").Value = -460.1428000000001
$ws.Range("N113").Value = -6940

# LTW row 2 (Leve Item ID 2631)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 1154924.9
$ws.Range("J2").Value = 1250918.5
$ws.Range("L2").Value = 1250918.5
$ws.Range("N2").Value = -1251142.5

# LTW row 22 (Leve Item ID 5277)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2011.8235
$ws.Range("I22").Value = 1699.7
$ws.Range("J22").Value = 2457.7144
$ws.Range("K22").Value = 1699.7
$ws.Range("L22").Value = 2457.7144
$ws.Range("M22").Value = -1404.7
$ws.Range("N22").Value = -3047.7144

# LTW row 27 (Leve Item ID 5277)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 2011.8235
$ws.Range("I27").Value = 1699.7
$ws.Range("J27").Value = 2457.7144
$ws.Range("K27").Value = 1699.7
$ws.Range("L27").Value = 2457.7144
$ws.Range("M27").Value = -1592.7
$ws.Range("N27").Value = -2671.7144

# LTW row 46 (Leve Item ID 5282)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1882.5
$ws.Range("I46").Value = 915
$ws.Range("J46").Value = 2850
$ws.Range("K46").Value = 915
$ws.Range("L46").Value = 2850
$ws.Range("M46").Value = -727
$ws.Range("N46").Value = -3226

# LTW row 75 (Leve Item ID 10833)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H75").Value = 43312.668
$ws.Range("J75").Value = 43312.668
$ws.Range("L75").Value = 43312.668
$ws.Range("N75").Value = -45184.668

# LTW row 78 (Leve Item ID 10833)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H78").Value = 43312.668
$ws.Range("J78").Value = 43312.668
$ws.Range("L78").Value = 129938.004
$ws.Range("N78").Value = -139298.004

# LTW row 122 (Leve Item ID 36247)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 7219.8
$ws.Range("I122").Value = 4799.5
$ws.Range("K122").Value = 14398.5
$ws.Range("M122").Value = -11948.5

# LTW row 132 (Leve Item ID 44058)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 5974.857
$ws.Range("I132").Value = 3207.5386
$ws.Range("J132").Value = 8373.200000000001
$ws.Range("K132").Value = 9622.6158
$ws.Range("L132").Value = 25119.6
$ws.Range("M132").Value = -7092.6158
$ws.Range("N132").Value = -30179.6

# WVR row 96 (Leve Item ID 19977)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 63990100
$ws.Range("I96").Value = 111112710
$ws.Range("J96").Value = 3403889.2
$ws.Range("K96").Value = 111112710
$ws.Range("L96").Value = 3403889.2
$ws.Range("M96").Value = -111111337
$ws.Range("N96").Value = -3406635.2

# WVR row 136 (Leve Item ID 44031)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 5102.4565
$ws.Range("I136").Value = 5047.107
$ws.Range("J136").Value = 5188.5557
$ws.Range("K136").Value = 15141.321
$ws.Range("L136").Value = 15565.6671
$ws.Range("M136").Value = -12591.321
$ws.Range("N136").Value = -20665.6671
